$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z3").Value = "a@gmail.com"
$ws.Range("Z3").Font.Color = 12611584
$ws.Range("Z3").Font.Underline = 2
Write-Host "done"
